# Update the Expenses sheet: replace the sample (Portuguese) expense rows
# with a new English expense dataset, and add a handful of new rows
# (8-11) plus a trailing "spacer" row (12) with a larger font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 2-7: overwrite date / amount / category / description ---
# Row 2
$ws.Range("A2").Value = 44972
$ws.Range("B2").Value = "Food"
$ws.Range("C2").Value = 45.67
$ws.Range("D2").Value = "Groceries"

# Row 3
$ws.Range("A3").Value = 45036
$ws.Range("B3").Value = "Shopping"
$ws.Range("C3").Value = 124.5
$ws.Range("D3").Value = "Clothing"

# Row 4
$ws.Range("A4").Value = 44990
$ws.Range("B4").Value = "Entertainment"
$ws.Range("C4").Value = 32.200000000000003
$ws.Range("D4").Value = "Movie Tickets"

# Row 5
$ws.Range("A5").Value = 45087
$ws.Range("B5").Value = "Utilities"
$ws.Range("C5").Value = 75
$ws.Range("D5").Value = "Electricity Bill"

# Row 6
$ws.Range("A6").Value = 45064
$ws.Range("B6").Value = "Travel"
$ws.Range("C6").Value = 230
$ws.Range("D6").Value = "Flight Tickets"

# Row 7
$ws.Range("A7").Value = 45010
$ws.Range("B7").Value = "Food"
$ws.Range("C7").Value = 18.75
$ws.Range("D7").Value = "Restaurant"

# --- New rows 8-11: copy the number formats from row 2 (date / currency), ---
# --- then fill in the values.                                             ---
$ws.Range("A2:D2").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A9:D9").PasteSpecial(-4122)
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Range("A11:D11").PasteSpecial(-4122)

# Row 8
$ws.Range("A8").Value = 44962
$ws.Range("B8").Value = "Health"
$ws.Range("C8").Value = 60
$ws.Range("D8").Value = "Doctor's Visit"

# Row 9
$ws.Range("A9").Value = 45119
$ws.Range("B9").Value = "Shopping"
$ws.Range("C9").Value = 89.88
$ws.Range("D9").Value = "Electronics"

# Row 10
$ws.Range("A10").Value = 44954
$ws.Range("B10").Value = "Entertainment"
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = "Concert Tickets"

# Row 11
$ws.Range("A11").Value = 45048
$ws.Range("B11").Value = "Travel"
$ws.Range("C11").Value = 320
$ws.Range("D11").Value = "Hotel Stay"

# --- Trailing spacer row 12: bigger font, taller row ---
$ws.Range("A12").Font.Size = 14
$ws.Rows.Item(12).RowHeight = 19

# Move the active selection to the new last cell, like the source workbook.
$null = $ws.Range("A12").Select()
